# ---------------------------------------------------------------------------
# Weekly NYPD CompStat report refresh: bump the report volume/week labels
# and replace the crime-complaint figures with the newly collected numbers.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels (rich-text shared strings) ------------------------------
# "Volume 32   Number  13" -> "Volume 32   Number  14"
$ws.Range("A8").Value = "Volume 32   Number  14"
# "Report Covering the Week  3/24/2025  Through  3/30/2025"
#  -> "Report Covering the Week  3/31/2025  Through  4/6/2025"
$ws.Range("C9").Value = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Crime complaint figures (rows 15-31), updated week-over-week ----------
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 51
$ws.Range("J15").Value = 42
$ws.Range("K15").Value = 21.428571428571
$ws.Range("L15").Value = 104
$ws.Range("M15").Value = 131.818181818182
$ws.Range("N15").Value = -7.272727272727
$ws.Range("C16").Value = 32
$ws.Range("D16").Value = 32
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 112
$ws.Range("G16").Value = 113
$ws.Range("H16").Value = -0.884955752212
$ws.Range("I16").Value = 363
$ws.Range("J16").Value = 382
$ws.Range("K16").Value = -4.973821989528
$ws.Range("L16").Value = -22.929936305732
$ws.Range("M16").Value = 23.050847457627
$ws.Range("N16").Value = -86.223908918406
$ws.Range("C17").Value = 54
$ws.Range("D17").Value = 40
$ws.Range("E17").Value = 35
$ws.Range("F17").Value = 183
$ws.Range("G17").Value = 167
$ws.Range("H17").Value = 9.580838323353
$ws.Range("I17").Value = 532
$ws.Range("J17").Value = 547
$ws.Range("K17").Value = -2.742230347349
$ws.Range("L17").Value = -4.488330341113
$ws.Range("M17").Value = 59.281437125748
$ws.Range("N17").Value = -34.076827757125
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 37
$ws.Range("E18").Value = -2.702702702702
$ws.Range("F18").Value = 143
$ws.Range("G18").Value = 127
$ws.Range("H18").Value = 12.59842519685
$ws.Range("I18").Value = 521
$ws.Range("J18").Value = 516
$ws.Range("K18").Value = 0.968992248062
$ws.Range("L18").Value = -6.126126126126
$ws.Range("M18").Value = -4.753199268738
$ws.Range("N18").Value = -83.31197950032
$ws.Range("C19").Value = 179
$ws.Range("D19").Value = 180
$ws.Range("E19").Value = -0.555555555555
$ws.Range("F19").Value = 741
$ws.Range("G19").Value = 708
$ws.Range("H19").Value = 4.661016949152
$ws.Range("I19").Value = 2531
$ws.Range("J19").Value = 2648
$ws.Range("K19").Value = -4.418429003021
$ws.Range("L19").Value = -12.452438602559
$ws.Range("M19").Value = -2.989651207359
$ws.Range("N19").Value = -68.578522656734
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -12.5
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 72
$ws.Range("J20").Value = 105
$ws.Range("K20").Value = -31.428571428571
$ws.Range("L20").Value = -41.935483870967
$ws.Range("M20").Value = -12.195121951219
$ws.Range("N20").Value = -95.593635250918
$ws.Range("C21").Value = 312
$ws.Range("D21").Value = 301
$ws.Range("E21").Value = 3.654485049833
$ws.Range("F21").Value = 1223
$ws.Range("G21").Value = 1168
$ws.Range("H21").Value = 4.708904109589
$ws.Range("I21").Value = 4073
$ws.Range("J21").Value = 4246
$ws.Range("K21").Value = -4.07442298634
$ws.Range("L21").Value = -11.954172070903
$ws.Range("M21").Value = 4.596815613764
$ws.Range("N21").Value = -75.067335945151
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = 17
$ws.Range("E22").Value = -41.176470588235
$ws.Range("F22").Value = 48
$ws.Range("G22").Value = 46
$ws.Range("H22").Value = 4.347826086956
$ws.Range("I22").Value = 165
$ws.Range("J22").Value = 174
$ws.Range("K22").Value = -5.172413793103
$ws.Range("L22").Value = -4.624277456647
$ws.Range("M22").Value = 11.486486486486
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = 50
$ws.Range("G23").Value = 28
$ws.Range("H23").Value = 35.714285714285
$ws.Range("I23").Value = 119
$ws.Range("J23").Value = 97
$ws.Range("K23").Value = 22.680412371134
$ws.Range("L23").Value = 22.680412371134
$ws.Range("M23").Value = 26.595744680851
$ws.Range("C24").Value = 396
$ws.Range("D24").Value = 357
$ws.Range("E24").Value = 10.924369747899
$ws.Range("F24").Value = 1599
$ws.Range("G24").Value = 1588
$ws.Range("H24").Value = 0.692695214105
$ws.Range("I24").Value = 5225
$ws.Range("J24").Value = 5527
$ws.Range("K24").Value = -5.46408539895
$ws.Range("L24").Value = 8.560149594847
$ws.Range("M24").Value = 29.267689262741
$ws.Range("C25").Value = 287
$ws.Range("D25").Value = 306
$ws.Range("E25").Value = -6.209150326797
$ws.Range("F25").Value = 1280
$ws.Range("G25").Value = 1376
$ws.Range("H25").Value = -6.976744186046
$ws.Range("I25").Value = 4282
$ws.Range("J25").Value = 4736
$ws.Range("K25").Value = -9.586148648648
$ws.Range("L25").Value = 4.84818805093
$ws.Range("C26").Value = 89
$ws.Range("D26").Value = 86
$ws.Range("E26").Value = 3.488372093023
$ws.Range("F26").Value = 432
$ws.Range("G26").Value = 367
$ws.Range("H26").Value = 17.711171662125
$ws.Range("I26").Value = 1262
$ws.Range("J26").Value = 1235
$ws.Range("K26").Value = 2.186234817813
$ws.Range("L26").Value = 6.497890295358
$ws.Range("M26").Value = 38.833883388338
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 6
$ws.Range("E27").Value = -16.666666666666
$ws.Range("F27").Value = 24
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 62
$ws.Range("J27").Value = 62
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 10.714285714285
$ws.Range("C28").Value = 17
$ws.Range("D28").Value = 21
$ws.Range("E28").Value = -19.047619047619
$ws.Range("F28").Value = 60
$ws.Range("G28").Value = 84
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 220
$ws.Range("J28").Value = 219
$ws.Range("K28").Value = 0.456621004566
$ws.Range("L28").Value = 0.456621004566
$ws.Range("D29").Value = 2
$ws.Range("G29").Value = 5
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -66.666666666666
$ws.Range("L29").Value = -72.727272727272
$ws.Range("M29").Value = -72.727272727272
$ws.Range("D30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = -62.5
$ws.Range("L30").Value = -66.666666666666
$ws.Range("M30").Value = -57.142857142857
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = -60
$ws.Range("F31").Value = 11
$ws.Range("G31").Value = 14
$ws.Range("H31").Value = -21.428571428571
$ws.Range("I31").Value = 32
$ws.Range("J31").Value = 36
$ws.Range("K31").Value = -11.111111111111
$ws.Range("L31").Value = -8.571428571428
$ws.Range("C33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("I33").Value = 5
$ws.Range("K33").Value = 25
$ws.Range("L33").Value = 66.666666666666

# --- Row 33 (Traffic Fatalities) also flips a couple of cells between ------
# numeric and text ("n/a"-style "***.*"/"0" placeholders) -------------------
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
